# Updated symbol list on Mon Jan  9 15:54:32 UTC 2023 with GitHub Actions
# Apply cell value updates to the cryptocurrency price/volume sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.70"
$ws.Range("E2").Value = "'5.04%"
$ws.Range("D3").Value = "'27.00"
$ws.Range("E3").Value = "'0.34%"
$ws.Range("D4").Value = "'4.937"
$ws.Range("E4").Value = "'5.29%"
$ws.Range("D5").Value = "'0.06381"
$ws.Range("E5").Value = "'4.27%"
$ws.Range("D6").Value = "'6.964"
$ws.Range("E6").Value = "'3.31%"
$ws.Range("D7").Value = "'3.349"
$ws.Range("E7").Value = "'5.88%"
$ws.Range("D8").Value = "'0.8857"
$ws.Range("E8").Value = "'4.21%"
$ws.Range("D9").Value = "'0.9514"
$ws.Range("E9").Value = "'4.65%"
$ws.Range("D10").Value = "'0.1473"
$ws.Range("E10").Value = "'4.03%"
$ws.Range("D11").Value = "'0.05149"
$ws.Range("E11").Value = "'6.78%"
$ws.Range("D12").Value = "'0.07407"
$ws.Range("E12").Value = "'4.53%"
$ws.Range("D13").Value = "'0.03119"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("D14").Value = "'0.09066"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("D15").Value = "'0.001552"
$ws.Range("E15").Value = "'1.11%"
$ws.Range("D16").Value = "'0.0006260"
$ws.Range("E16").Value = "'1.44%"
$ws.Range("D17").Value = "'0.005978"
$ws.Range("E17").Value = "'0.23%"
$ws.Range("D18").Value = "'3.500"
$ws.Range("E18").Value = "'1.44%"
$ws.Range("E19").Value = "'7.03%"
$ws.Range("D21").Value = "'0.1328"
$ws.Range("E21").Value = "'3.76%"
$ws.Range("D22").Value = "'3.867"
$ws.Range("E22").Value = "'-6.05%"
$ws.Range("D23").Value = "'0.04346"
$ws.Range("E23").Value = "'2.49%"
$ws.Range("E24").Value = "'-0.14%"
$ws.Range("D25").Value = "'0.003658"
$ws.Range("E25").Value = "'-10.35%"
$ws.Range("D26").Value = "'0.0001196"
$ws.Range("E26").Value = "'-0.13%"
$ws.Range("D27").Value = "'0.0001689"
$ws.Range("E27").Value = "'0.67%"
$ws.Range("D40").Value = "'0.04076"
$ws.Range("E40").Value = "'3.45%"
$ws.Range("D41").Value = "'0.006616"
$ws.Range("E41").Value = "'58.70%"
$ws.Range("D42").Value = "'0.1175"
$ws.Range("E42").Value = "'4.86%"
$ws.Range("D43").Value = "'0.002343"
$ws.Range("D44").Value = "'0.01254"
$ws.Range("E44").Value = "'7.06%"
$ws.Range("D45").Value = "'0.00005243"
$ws.Range("E45").Value = "'3.08%"
$ws.Range("E46").Value = "'0.19%"
$ws.Range("D47").Value = "'2.352"
$ws.Range("E47").Value = "'835.11%"
$ws.Range("D48").Value = "'0.02254"
$ws.Range("E48").Value = "'-7.86%"
$ws.Range("E49").Value = "'0.19%"
$ws.Range("E50").Value = "'0.11%"
